$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph and insert a new
# ListBullet paragraph right after it, listing the two instructors
# (each on its own line, separated by a manual line break), mirroring the
# existing "Requisitos"/"Avaliação" bullet-list paragraphs already in the
# document.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Docente(s) Responsável(eis)*") {
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Style = "ListBullet"

        $newPara.Range.Select()
        $sel = $word.Selection
        $sel.TypeText("5983729 - Fernando Vernilli Junior")
        $sel.InsertBreak(6)
        $sel.TypeText("1922320 - Sebastiao Ribeiro")

        break
    }
}
